$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Move Robot21 to location (6, 12) and remove the toolkit.'
$ws.Range("B2").Value = 96.392025
$ws.Range("C2").Value = 13546
$ws.Range("D2").Value = "'0.04161"
$ws.Range("E2").Value = 'a41b30bc-f7db-45b9-8182-a512b305d635'

$ws.Range("A3").Value = 'Move Robot41 to location (1, 11) and remove the liquid spill.'
$ws.Range("B3").Value = 107.572532
$ws.Range("C3").Value = 13620
$ws.Range("D3").Value = "'0.04257"
$ws.Range("E3").Value = '1f84bbbf-e00a-45ec-80e7-1f6ed186be91'

$ws.Range("A4").Value = 'Move Robot9 to location (1, 5) and remove the large debris.'
$ws.Range("B4").Value = 66.73797
$ws.Range("C4").Value = 10479
$ws.Range("D4").Value = "'0.03246"
$ws.Range("E4").Value = '6bcf9400-440f-40d2-b339-cda51dd30392'

$ws.Range("A5").Value = 'Move Robot42 to location (1, 11) and remove the dust.'
$ws.Range("B5").Value = 96.28057
$ws.Range("C5").Value = 13608
$ws.Range("D5").Value = "'0.04272"
$ws.Range("E5").Value = 'ba6f88ce-7e0e-465a-8ab9-f6e0ac7b52c1'

$ws.Range("A6").Value = 'Move Robot32 to location (3, 9) and remove the grass.'
$ws.Range("B6").Value = 98.646495
$ws.Range("C6").Value = 13731
$ws.Range("D6").Value = "'0.04476"
$ws.Range("E6").Value = 'a76c3491-e4ff-43a7-8b16-29c17b37d81c'

$ws.Range("A7").Value = 'Move Robot14 to location (11, 12) and remove the small debris.'
$ws.Range("B7").Value = 98.8925
$ws.Range("C7").Value = 14149
$ws.Range("D7").Value = "'0.04533"
$ws.Range("E7").Value = '8375493b-21ac-4113-84ee-4f6c15df3ce0'

$ws.Range("A8").Value = 'Move Robot39 to location (6, 4) and remove the vehicle.'
$ws.Range("B8").Value = 97.838499
$ws.Range("C8").Value = 13940
$ws.Range("D8").Value = "'0.04479"
$ws.Range("E8").Value = '1aa344cc-78cd-4399-8244-e3e0ce1de67f'

$ws.Range("A9").Value = 'Move Robot15 to location (11, 2) and remove the construction materials.'
$ws.Range("B9").Value = 88.1175
$ws.Range("C9").Value = 13984
$ws.Range("D9").Value = "'0.04434"
$ws.Range("E9").Value = '25d7bd10-7160-4c30-b979-9a59accf2fa9'

$ws.Range("A10").Value = 'Move Robot2 to location (1, 10) and remove the tree branches.'
$ws.Range("B10").Value = 101.209372
$ws.Range("C10").Value = 14049
$ws.Range("D10").Value = "'0.04518"
$ws.Range("E10").Value = 'd67e2002-a34f-4197-8b17-2c6d52e6cbc7'

$ws.Range("A11").Value = 'Move Robot26 to location (1, 3) and remove the screws.'
$ws.Range("B11").Value = 99.200013
$ws.Range("C11").Value = 13552
$ws.Range("D11").Value = "'0.04236"
$ws.Range("E11").Value = '9707ff1c-3274-4a14-9e91-8919b53217c3'

Write-Host "Applied batch 6 updates"